$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6) entirely - the dataset is trimmed by one row.
$ws.Rows.Item(6).Delete()

# Apply "custom accuracy" - round row 5's measurement values down to 2 decimal places
# (previously stored with 3 decimals).
$ws.Range("B5").Value = 20.18
$ws.Range("C5").Value = 14.82
$ws.Range("D5").Value = 1.22
$ws.Range("E5").Value = 43.86
$ws.Range("F5").Value = 35.69
$ws.Range("G5").Value = 15.81
$ws.Range("H5").Value = 61.91
$ws.Range("I5").Value = 24.43
$ws.Range("J5").Value = 10.82
$ws.Range("K5").Value = 15.97
$ws.Range("L5").Value = 17.59
$ws.Range("M5").Value = 18.6
$ws.Range("N5").Value = 5.07
$ws.Range("O5").Value = 15.79
$ws.Range("P5").Value = 22.44
$ws.Range("Q5").Value = 13.36
$ws.Range("R5").Value = 0.82
$ws.Range("S5").Value = 0.83
$ws.Range("T5").Value = 233.04
$ws.Range("U5").Value = 44.18
$ws.Range("V5").Value = 14.58
$ws.Range("W5").Value = 29.63
$ws.Range("X5").Value = 15.51
$ws.Range("Y5").Value = 2.41
$ws.Range("Z5").Value = 30.09
$ws.Range("AA5").Value = 12.87
$ws.Range("AB5").Value = 11.43
$ws.Range("AC5").Value = 13.45
$ws.Range("AD5").Value = 18.49
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 56.42
$ws.Range("AG5").Value = 8.17
$ws.Range("AH5").Value = 18.22

# Column AB (28th column) narrows from 8 to 7 characters - match the width already
# used by the other "7"-wide columns (e.g. column D) instead of a hard-coded magic
# number, since ColumnWidth uses a different unit than the stored sheet width.
$ws.Columns.Item(28).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
